# Updates the cryptos price/volume snapshot (and two ranking swaps) to the
# latest scrape, per commit "Updated cryptos list on Fri Nov 15 13:37:53 UTC 2024
# with GitHub Actions".
#
# Numeric-looking Price values are prefixed with a leading apostrophe so Excel
# keeps them as text (matching the source workbook, where these cells are
# plain text, not numbers -- e.g. "1.00" / "214.39"), instead of silently
# coercing them to numbers and losing the trailing zeros / formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "89.791.66"
$ws.Range("E2").Value = "  -1.17%  "
$ws.Range("D3").Value = "3.100.35"
$ws.Range("E3").Value = "  -2.48%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'214.39"
$ws.Range("E5").Value = "  -0.82%  "
$ws.Range("D6").Value = "'619.50"
$ws.Range("E6").Value = "  -1.89%  "
$ws.Range("E7").Value = "  -5.92%  "
$ws.Range("D8").Value = "'0.873"
$ws.Range("E8").Value = "  +20.37%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").Value = "3.099.66"
$ws.Range("E10").Value = "  -2.35%  "
$ws.Range("D11").Value = "'0.665"
$ws.Range("E11").Value = "  +17.19%  "
$ws.Range("E12").Value = "  +3.28%  "
$ws.Range("E13").Value = "  -5.10%  "
$ws.Range("D14").Value = "'5.39"
$ws.Range("E14").Value = "  +1.09%  "
$ws.Range("D15").Value = "89.874.51"
$ws.Range("D16").Value = "'32.60"
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("D17").Value = "3.680.07"
$ws.Range("E17").Value = "  -2.59%  "
$ws.Range("D18").Value = "3.093.34"
$ws.Range("E18").Value = "  -3.29%  "
$ws.Range("D19").Value = "'3.41"
$ws.Range("E19").Value = "  +1.93%  "
$ws.Range("D20").Value = "'0.0000215"
$ws.Range("E20").Value = "  +1.45%  "
$ws.Range("D21").Value = "'13.57"
$ws.Range("E21").Value = "  +1.06%  "
$ws.Range("D22").Value = "'433.30"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").Value = "'8.29"
$ws.Range("E23").Value = "  -2.03%  "
$ws.Range("D24").Value = "'5.02"
$ws.Range("E24").Value = "  +0.61%  "
$ws.Range("D25").Value = "'5.56"
$ws.Range("E25").Value = "  +4.82%  "
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").Value = "'86.15"
$ws.Range("E26").Value = "  +6.72%  "
$ws.Range("B27").Value = "Aptos"
$ws.Range("C27").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D27").Value = "'12.19"
$ws.Range("E27").Value = "  +4.66%  "
$ws.Range("D28").Value = "3.284.33"
$ws.Range("E28").Value = "  -1.86%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("D30").Value = "'1.10"
$ws.Range("E30").Value = "  +9.43%  "
$ws.Range("D31").Value = "'0.164"
$ws.Range("E31").Value = "  +3.04%  "
$ws.Range("D32").Value = "'8.17"
$ws.Range("E32").Value = "  -2.00%  "
$ws.Range("D33").Value = "'515.74"
$ws.Range("E33").Value = "  +0.70%  "
$ws.Range("E34").Value = "  -8.04%  "
$ws.Range("D35").Value = "'6.76"
$ws.Range("E35").Value = "  -3.02%  "
$ws.Range("D36").Value = "'23.03"
$ws.Range("E36").Value = "  +3.40%  "
$ws.Range("B37").Value = "PancakeSwap"
$ws.Range("C37").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D37").Value = "'1.81"
$ws.Range("E37").Value = "  -4.01%  "
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").Value = "'1.25"
$ws.Range("E38").Value = "  -2.65%  "
$ws.Range("E39").Value = "  +4.50%  "
$ws.Range("D40").Value = "'22.30"
$ws.Range("E40").Value = "  -0.30%  "
$ws.Range("E41").Value = "  +0.17%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("E43").Value = "  +13.97%  "
$ws.Range("D44").Value = "'0.373"
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("E45").Value = "  -3.07%  "
$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").Value = "'146.06"
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").Value = "'0.0707"
$ws.Range("E47").Value = "  +14.56%  "
$ws.Range("D48").Value = "'43.57"
$ws.Range("E48").Value = "  -0.86%  "
$ws.Range("E49").Value = "  +2.18%  "
$ws.Range("D50").Value = "'159.82"
$ws.Range("E50").Value = "  -5.39%  "
$ws.Range("E51").Value = "  +0.43%  "
